$d = $word.ActiveDocument

# 1) "hopi mayor/kuidadó no ta kustumá" -> "hopi mayor/dunadónan di kuido no ta kustumá"
$d.Content.Find.Execute("mayor/kuidadó no ta kustumá", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "mayor/dunadónan di kuido no ta kustumá", 2)

# 2) "mayornan/kuidadónan tambe por permití" -> "mayornan/dunadónan di kuido tambe por permití"
$d.Content.Find.Execute("mayornan/kuidadónan tambe por permití", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "mayornan/dunadónan di kuido tambe por permití", 2)

# 3) "pa nan por kontestá na nan yunan na un manera" -> "pa nan por kontestá nan yunan na un manera"
$d.Content.Find.Execute("kontestá na nan yunan na un manera", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "kontestá nan yunan na un manera", 2)

# 4a) "no ta kustumbrá ku ta puntra" -> "no ta kustumá ku ta puntra"
$d.Content.Find.Execute("kustumbrá ku ta puntra", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "kustumá ku ta puntra", 2)

# 4b) "ku hende ta skuchá nan" -> "ku hende ta skucha nan"
$d.Content.Find.Execute("ku hende ta skuchá nan", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "ku hende ta skucha nan", 2)
